$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: date 2021-02-22, "Made a Player Page for Spotify Authentication [template only]"
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)
$ws.Cells.Item(9, 1).Value = 44249
$ws.Cells.Item(9, 2).Value = "Made a Player Page for Spotify Authentication [template only]"

# Row 10: "Read about React Context API", 0.75 hours
$ws.Cells.Item(10, 2).Value = "Read about React Context API"
$ws.Cells.Item(10, 3).Value = 0.75

# Row 11: date 2021-02-23, "Used Context API instead of regular States"
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = 44250
$ws.Cells.Item(11, 2).Value = "Used Context API instead of regular States"

# Row 12: "Passed Spotify as a prop to the Player Component"
$ws.Cells.Item(12, 2).Value = "Passed Spotify as a prop to the Player Component"

# Row 13: "Bypassed the login page by setting the initial value of Token (this step has to be undone later on)", 2 hours
$ws.Cells.Item(13, 2).Value = "Bypassed the login page by setting the initial value of Token (this step has to be undone later on)"
$ws.Cells.Item(13, 3).Value = 2
$ws.Rows.Item(13).RowHeight = 28.8

# Row 14: "To Decide" / "Create Player, SideBar, SongRow"
# (shared string for "Create Player, SideBar, SongRow" must be registered before "To Decide")
$ws.Cells.Item(14, 2).Value = "Create Player, SideBar, SongRow"
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(14, 1).PasteSpecial(-4122)
$ws.Cells.Item(14, 1).Value = "To Decide"

# Row 15: "To Decide" / "Host it on Firebase"
$ws.Cells.Item(15, 1).Value = "To Decide"
$ws.Cells.Item(15, 2).Value = "Host it on Firebase"

# Row 16: "To Decide" / "Share your Spotify Clone App and Knowlegde on Medium"
$ws.Cells.Item(16, 1).Value = "To Decide"
$ws.Cells.Item(16, 2).Value = "Share your Spotify Clone App and Knowlegde on Medium"

$ws.Range("B16").Select()
